$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.763564666666667
$ws.Range("H2").Value = 8.290694
$ws.Range("I2").Value = 0.009909756414635561
$ws.Range("J2").Value = 0.009909756414635559
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 213.0418828541756
$ws.Range("R2").Value = 1917.37694568758
$ws.Range("S2").Value = 0.002382123770207231
$ws.Range("T2").Value = 0.002382123770207231
$ws.Range("G3").Value = 2.763564666666667
$ws.Range("H3").Value = 8.290694
$ws.Range("I3").Value = 0.009909756414635561
$ws.Range("J3").Value = 0.009909756414635559
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 280.7230020130809
$ws.Range("R3").Value = 2526.507018117728
$ws.Range("S3").Value = 0.003138898919688109
$ws.Range("T3").Value = 0.003138898919688108
$ws.Range("G4").Value = 2.763564666666667
$ws.Range("H4").Value = 8.290694
$ws.Range("I4").Value = 0.009909756414635561
$ws.Range("J4").Value = 0.009909756414635559
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 392.5002167217103
$ws.Range("R4").Value = 3532.501950495392
$ws.Range("S4").Value = 0.004388733724740221
$ws.Range("T4").Value = 0.00438873372474022
$ws.Range("G5").Value = 266.1315866666666
$ws.Range("H5").Value = 798.3947599999999
$ws.Range("I5").Value = 0.9543106517164206
$ws.Range("J5").Value = 0.9543106517164204
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 20515.95716007702
$ws.Range("R5").Value = 184643.6144406932
$ws.Range("S5").Value = 0.2293987856511044
$ws.Range("T5").Value = 0.2293987856511043
$ws.Range("G6").Value = 266.1315866666666
$ws.Range("H6").Value = 798.3947599999999
$ws.Range("I6").Value = 0.9543106517164206
$ws.Range("J6").Value = 0.9543106517164204
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("Q6").Value = 27033.65650917923
$ws.Range("R6").Value = 243302.9085826131
$ws.Range("S6").Value = 0.3022763172357643
$ws.Range("T6").Value = 0.3022763172357641
$ws.Range("G7").Value = 266.1315866666666
$ws.Range("H7").Value = 798.3947599999999
$ws.Range("I7").Value = 0.9543106517164206
$ws.Range("J7").Value = 0.9543106517164204
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("Q7").Value = 37797.81479445241
$ws.Range("S7").Value = 0.422635548829552
$ws.Range("T7").Value = 0.422635548829552
$ws.Range("I8").Value = 0.03577959186894402
$ws.Range("J8").Value = 0.03577959186894401
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 769.1966684727144
$ws.Range("R8").Value = 6922.770016254431
$ws.Range("S8").Value = 0.008600757951370833
$ws.Range("T8").Value = 0.008600757951370832
$ws.Range("I9").Value = 0.03577959186894402
$ws.Range("J9").Value = 0.03577959186894401
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.01133312642260744
$ws.Range("T9").Value = 0.01133312642260744
$ws.Range("I10").Value = 0.03577959186894402
$ws.Range("J10").Value = 0.03577959186894401
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.01584570749496575
$ws.Range("T10").Value = 0.01584570749496575
